# Remove the "common_name" column (column D) from every worksheet.
# Deleting the entire column shifts tot_fronds (previously column E)
# left into column D, matching the diff.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Columns.Item(4).Delete()
}
